$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update specific "subcategory" (column H) values
$ws.Range("H4").Value = "photo(s)"
$ws.Range("H5").Value = "drawing(s)"
$ws.Range("H14").Value = "line graph(s)"
$ws.Range("H15").Value = "line graph(s)"
$ws.Range("H16").Value = "scatter plot(s)"
$ws.Range("H17").Value = "line graph(s)"
$ws.Range("H18").Value = "line graph(s)"
$ws.Range("H20").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H26").Value = "drawing(s)"
$ws.Range("H27").Value = "drawing(s)"
$ws.Range("H28").Value = "photo(s)"
$ws.Range("H34").Value = "drawing(s)"
$ws.Range("H43").Value = "drawing(s)"
$ws.Range("H45").Value = "drawing(s)"
$ws.Range("H46").Value = "drawing(s)"
$ws.Range("H47").Value = "drawing(s)"
$ws.Range("H60").Value = "line graph(s)"
$ws.Range("H61").Value = "line graph(s)"
$ws.Range("H66").Value = "data display"
$ws.Range("H71").Value = "line graph(s)"
$ws.Range("H75").Value = "line graph(s)"

# Remove entire column I ("is_viewed") -- shifts nothing else since it's the last column
$ws.Range("I1:I83").Delete()
